# Rebuild the LOM3036 worksheet content to match the revised course-plan data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate: clear all cell content/formatting and any custom row heights.
$ws.Cells.Clear()
for ($r = 1; $r -le 26; $r++) { $ws.Rows.Item($r).AutoFit() }

# Helper: write a value as plain text (never let Excel auto-convert to number/date)
# and land on a specific look: Bold / WrapText / font color, vertical-top aligned.
function Set-CellText($rng, $text, $bold, $wrap, $colorIndex) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
    $rng.Font.Bold = $bold
    $rng.WrapText = $wrap
    if ($colorIndex -eq $null) {
        $rng.Font.ColorIndex = -4105
    } else {
        $rng.Font.Color = $colorIndex
    }
    $rng.VerticalAlignment = -4160
}

function Set-ColA($rng, $text) { Set-CellText $rng $text $true  $false $null }
function Set-ColB($rng, $text) { Set-CellText $rng $text $false $true  $null }
function Set-ColC($rng, $text) { Set-CellText $rng $text $false $true  255 }

# Row 1
Set-ColA $ws.Range("B1") "Ementa atual:"
Set-ColA $ws.Range("C1") "Ementa modificada (dados modificados em vermelho):"

# Row 2
Set-ColB $ws.Range("B2") "LOM3036"
Set-ColC $ws.Range("C2") "LOM3036"

# Row 3
Set-ColA $ws.Range("A3") "Nome:"
Set-ColB $ws.Range("B3") " Propriedades Mecânicas"
Set-ColC $ws.Range("C3") " Propriedades Mecânicas"

# Row 4
Set-ColA $ws.Range("A4") "Name:"
Set-ColB $ws.Range("B4") "Mechanical Properties"
Set-ColC $ws.Range("C4") "Mechanical Properties"

# Row 5
Set-ColA $ws.Range("A5") "Créditos-aula:"
Set-ColB $ws.Range("B5") "4"
Set-ColC $ws.Range("C5") "4"

# Row 6
Set-ColA $ws.Range("A6") "Créditos-trabalho"
Set-ColB $ws.Range("B6") "0"
Set-ColC $ws.Range("C6") "0"

# Row 7
Set-ColA $ws.Range("A7") "Carga horária:"
Set-ColB $ws.Range("B7") "60 h"
Set-ColC $ws.Range("C7") "60 h"

# Row 8
Set-ColA $ws.Range("A8") "Ativação:"
Set-ColB $ws.Range("B8") "01/01/2018"
Set-ColC $ws.Range("C8") "01/01/2018"

# Row 9
Set-ColA $ws.Range("A9") "Semestre ideal:"
Set-ColB $ws.Range("B9") "EM-6"
Set-ColC $ws.Range("C9") "EM-6"

# Row 10
Set-ColA $ws.Range("A10") "Objetivos:"
Set-ColB $ws.Range("B10") "Apresentar os princípios básicos de Ciências dos Materiais, destacando a correlação entre o comportamento mecânico dos metais e os aspectos microestruturais, para aplicação em Engenharia."
Set-ColC $ws.Range("C10") "Apresentar os princípios básicos de Ciências dos Materiais, destacando a correlação entre o comportamento mecânico dos metais e os aspectos microestruturais, para aplicação em Engenharia."
$ws.Rows.Item(10).RowHeight = 60

# Row 11
Set-ColA $ws.Range("A11") "Objectives:"
$ws.Rows.Item(11).RowHeight = 60

# Row 12
Set-ColA $ws.Range("A12") "Docentes responsáveis:"

# Row 13
Set-ColB $ws.Range("B13") "7459752 - Maria Ismenia Sodero Toledo Faria"
Set-ColC $ws.Range("C13") "7459752 - Maria Ismenia Sodero Toledo Faria"

# Row 14
Set-ColB $ws.Range("B14") "5840622 - Miguel Justino Ribeiro Barboza"
Set-ColC $ws.Range("C14") "5840622 - Miguel Justino Ribeiro Barboza"

# Row 15
Set-ColA $ws.Range("A15") "Programa resumido:"
Set-ColB $ws.Range("B15") "1. Introdução ao conceito de propriedades mecânicas.2. Deformação plástica de monocristais e policristais.3. Teoria das discordâncias.4. Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos6. Influências ambientais e térmicas no comportamento mecânico. Análise de falhas."
Set-ColC $ws.Range("C15") "1. Introdução ao conceito de propriedades mecânicas.2. Deformação plástica de monocristais e policristais.3. Teoria das discordâncias.4. Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos6. Influências ambientais e térmicas no comportamento mecânico. Análise de falhas."
$ws.Rows.Item(15).RowHeight = 60

# Row 16
Set-ColA $ws.Range("A16") "Short syllabus:"
$ws.Rows.Item(16).RowHeight = 60

# Row 17
Set-ColA $ws.Range("A17") "Programa:"
Set-ColB $ws.Range("B17") "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas. Comportamento elástico e plástico de metais e ligas. Relações entre tensão e deformação uniaxiais para regime plástico.2.DEFORMAÇÃO PLÁSTICA DE MONOCRISTAIS E POLICRISTAIS: Deformação plástica e encruamento de monocristais. Sistemas de deslizamento. Deformação por maclação e movimentação de discordâncias. Movimento relativo de grãos.3.TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4.MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão e à formação de células e subgrãos. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Diagrama Ferro-Carbono. Curvas TTT. Aços comuns e especiais. Tratamentos térmicos em aços; Transformação martensítica.5.COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil.6.Influências ambientais e térmicas sobre o comportamento mecânico dos metais. Análise de falhas em componentes."
Set-ColC $ws.Range("C17") "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas. Comportamento elástico e plástico de metais e ligas. Relações entre tensão e deformação uniaxiais para regime plástico.2.DEFORMAÇÃO PLÁSTICA DE MONOCRISTAIS E POLICRISTAIS: Deformação plástica e encruamento de monocristais. Sistemas de deslizamento. Deformação por maclação e movimentação de discordâncias. Movimento relativo de grãos.3.TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4.MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão e à formação de células e subgrãos. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Diagrama Ferro-Carbono. Curvas TTT. Aços comuns e especiais. Tratamentos térmicos em aços; Transformação martensítica.5.COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil.6.Influências ambientais e térmicas sobre o comportamento mecânico dos metais. Análise de falhas em componentes."
$ws.Rows.Item(17).RowHeight = 120

# Row 18
Set-ColA $ws.Range("A18") "Syllabus:"
$ws.Rows.Item(18).RowHeight = 120

# Row 19
Set-ColA $ws.Range("A19") "Avaliação:"

# Row 20
Set-ColA $ws.Range("A20") "Método:"
Set-ColB $ws.Range("B20") "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
Set-ColC $ws.Range("C20") "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."

# Row 21
Set-ColA $ws.Range("A21") "Critério:"
Set-ColB $ws.Range("B21") "A média do semestre será computada com base na relação:M=(P1+2P2)/3"
Set-ColC $ws.Range("C21") "A média do semestre será computada com base na relação:M=(P1+2P2)/3"
$ws.Rows.Item(21).RowHeight = 60

# Row 22
Set-ColA $ws.Range("A22") "Norma de recuperação:"
Set-ColB $ws.Range("B22") "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2"
Set-ColC $ws.Range("C22") "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2"
$ws.Rows.Item(22).RowHeight = 60

# Row 23
Set-ColA $ws.Range("A23") "Bibliografia:"
Set-ColB $ws.Range("B23") "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009.2. Hearn, E.J. Mechanics of Materials: An Introduction to the Mechanics of Elastic and Plastic Deformation of Solids and Structural Components, Pergamon Press, 1985.3. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1981.4. Hull, D. Introduction to Dislocations, Pergamon Press, 1965.5. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967.6. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.7. Van Vlack, L.H. Princípios de Ciência dos materiais, Ed. Edgard Blucher Ltda., 1970.8. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008.9. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall,1988.10. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008.11. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
Set-ColC $ws.Range("C23") "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009.2. Hearn, E.J. Mechanics of Materials: An Introduction to the Mechanics of Elastic and Plastic Deformation of Solids and Structural Components, Pergamon Press, 1985.3. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1981.4. Hull, D. Introduction to Dislocations, Pergamon Press, 1965.5. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967.6. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.7. Van Vlack, L.H. Princípios de Ciência dos materiais, Ed. Edgard Blucher Ltda., 1970.8. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008.9. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall,1988.10. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008.11. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$ws.Rows.Item(23).RowHeight = 120

# Row 24
Set-ColA $ws.Range("A24") "Requisitos:"

# Row 25
Set-ColB $ws.Range("B25") "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
Set-ColC $ws.Range("C25") "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Rows.Item(25).RowHeight = 30

# Row 26
Set-ColB $ws.Range("B26") "LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)`n"
Set-ColC $ws.Range("C26") "LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)`n"
$ws.Rows.Item(26).RowHeight = 30
